$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and hourly volume-change (E) columns.
# D values are prefixed with a leading apostrophe (the standard Excel
# text-entry marker) to force text storage, matching the workbook's
# original data, which is not auto-converted to a number even for
# numeric-looking strings such as 279.70 or 30.257.06.

$ws.Range("D2").Value = "'30.257.06"
$ws.Range("E2").Value = "  +0.47%  "

$ws.Range("D3").Value = "'1.862.54"
$ws.Range("E3").Value = "  +0.39%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'235.62"
$ws.Range("E5").Value = "  +0.83%  "

$ws.Range("E6").Value = "  -0.02%  "

$ws.Range("D7").Value = "'0.4671"
$ws.Range("E7").Value = "  -0.50%  "

$ws.Range("D8").Value = "'0.2839"
$ws.Range("E8").Value = "  +1.37%  "

$ws.Range("D9").Value = "'0.06525"

$ws.Range("D10").Value = "'21.67"
$ws.Range("E10").Value = "  +8.95%  "

$ws.Range("D11").Value = "'0.07940"
$ws.Range("E11").Value = "  +1.94%  "

$ws.Range("D12").Value = "'97.33"
$ws.Range("E12").Value = "  +0.69%  "

$ws.Range("D13").Value = "'1.871.27"
$ws.Range("E13").Value = "  +0.83%  "

$ws.Range("D14").Value = "'5.148"
$ws.Range("E14").Value = "  +1.23%  "

$ws.Range("D15").Value = "'0.6786"
$ws.Range("E15").Value = "  +2.33%  "

$ws.Range("D16").Value = "'279.70"
$ws.Range("E16").Value = "  -0.44%  "

$ws.Range("D17").Value = "'30.256.35"
$ws.Range("E17").Value = "  +0.35%  "

$ws.Range("D18").Value = "'13.46"
$ws.Range("E18").Value = "  +7.44%  "

$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "  +0.01%  "

$ws.Range("D20").Value = "'5.389"
$ws.Range("E20").Value = "  -1.48%  "

$ws.Range("D21").Value = "'2.113.98"
$ws.Range("E21").Value = "  +0.85%  "

$ws.Range("D22").Value = "'0.000007305"
$ws.Range("E22").Value = "  +1.28%  "

$ws.Range("E23").Value = "  -0.04%  "

$ws.Range("D24").Value = "'6.172"
$ws.Range("E24").Value = "  +0.85%  "

$ws.Range("D25").Value = "'167.29"
$ws.Range("E25").Value = "  -0.32%  "

$ws.Range("D26").Value = "'9.178"
$ws.Range("E26").Value = "  -0.98%  "

$ws.Range("E27").Value = "  +0.91%  "

$ws.Range("D28").Value = "'1.929"
$ws.Range("E28").Value = "  +0.99%  "

$ws.Range("D29").Value = "'1.389"
$ws.Range("E29").Value = "  +3.55%  "

$ws.Range("D30").Value = "'0.09728"
$ws.Range("E30").Value = "  +1.55%  "

$ws.Range("D31").Value = "'4.369"
$ws.Range("E31").Value = "  -0.78%  "

$ws.Range("D32").Value = "'1.476"
$ws.Range("E32").Value = "  +0.60%  "

$ws.Range("D33").Value = "'4.060"
$ws.Range("E33").Value = "  -0.53%  "

$ws.Range("D34").Value = "'0.04740"
$ws.Range("E34").Value = "  +1.75%  "

$ws.Range("D35").Value = "'1.129"
$ws.Range("E35").Value = "  +3.56%  "

$ws.Range("D36").Value = "'0.7059"
$ws.Range("E36").Value = "  +1.90%  "

$ws.Range("E37").Value = "  -0.03%  "

$ws.Range("D38").Value = "'0.01866"
$ws.Range("E38").Value = "  +1.20%  "

$ws.Range("D39").Value = "'2.583"
$ws.Range("E39").Value = "  +3.14%  "

$ws.Range("D40").Value = "'6.303"
$ws.Range("E40").Value = "  +0.61%  "

$ws.Range("D41").Value = "'75.06"
$ws.Range("E41").Value = "  +4.88%  "

$ws.Range("E42").Value = "  +1.27%  "

$ws.Range("D43").Value = "'0.8508"
$ws.Range("E43").Value = "  -0.46%  "

$ws.Range("D44").Value = "'0.4181"
$ws.Range("E44").Value = "  +1.07%  "

$ws.Range("E45").Value = "  -0.02%  "

$ws.Range("D46").Value = "'103.50"
$ws.Range("E46").Value = "  -0.48%  "

$ws.Range("D47").Value = "'964.66"
$ws.Range("E47").Value = "  -5.01%  "

$ws.Range("D48").Value = "'7.170"
$ws.Range("E48").Value = "  +0.01%  "

$ws.Range("D49").Value = "'9.295"
$ws.Range("E49").Value = "  +4.51%  "

$ws.Range("D50").Value = "'34.05"
$ws.Range("E50").Value = "  +1.35%  "

$ws.Range("D51").Value = "'0.1131"
$ws.Range("E51").Value = "  -0.29%  "
